$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the four "(26)" sheets, prefixing them with "!" (now
#    treated like the workbook's other archived/inactive sheets).
# ------------------------------------------------------------------
$ws65  = $wb.Worksheets.Item("S_65C (26)")
$ws90  = $wb.Worksheets.Item("S_90C (26)")
$ws110 = $wb.Worksheets.Item("S_110C (26)")
$ws130 = $wb.Worksheets.Item("S_130C (26)")

$ws65.Name  = "!S_65C (26)"
$ws90.Name  = "!S_90C (26)"
$ws110.Name = "!S_110C (26)"
$ws130.Name = "!S_130C (26)"

# ------------------------------------------------------------------
# 2. Update the "primary ref data check?" flag cell on each of the
#    four sheets: it used to point at the "primary ref data check?"
#    shared string with the green "Good" style, now it should read
#    "Same as Hilic 2001 (20)" with the yellow "Neutral" style.
# ------------------------------------------------------------------
$ws65.Range("F4").Value = "Same as Hilic 2001 (20)"
$ws65.Range("F4").Style = "Neutral"

$ws90.Range("G4").Value = "Same as Hilic 2001 (20)"
$ws90.Range("G4").Style = "Neutral"

$ws110.Range("G4").Value = "Same as Hilic 2001 (20)"
$ws110.Range("G4").Style = "Neutral"

$ws130.Range("G4").Value = "Same as Hilic 2001 (20)"
$ws130.Range("G4").Style = "Neutral"

# ------------------------------------------------------------------
# 3. Update the per-sheet selections (and, as a side effect of
#    selecting a range on a sheet, which sheet tab is active).
#    Order matters: the last sheet selected ends up as the active
#    tab, so select !S_130C (26) last.
# ------------------------------------------------------------------
$ws65.Range("F4").Select()
$ws90.Range("G4").Select()
$ws110.Range("G4").Select()
$ws130.Range("G4").Select()
$ws130.Range("L10").Select()
